# sales-invoice.xlsx: refresh invoice date/number, customer name, address,
# contact number, and the single line-item (name/qty/unit-price/total),
# per the corrected order details.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header: invoice date + invoice number
$ws.Range("F2").Value = "2021-06-10 "
$ws.Range("F3").Value = "INV202106101423"

# BILL TO block
$ws.Range("A9").Value = "Sant Anurag Deo"
$ws.Range("A11").Value = "Bangalore East,Karnataka"
$ws.Range("A13").Value = "'9900019861"

# Salesperson / customer row
$ws.Range("B16").Value = "Sant Anurag Deo"
$ws.Range("D16").Value = "'9900019861"

# Line item 1: description, qty stays, unit price, tax stays, total
$ws.Range("B19").Value = "Introduction To Vihangam Yoga"
$ws.Range("D19").Value = "'50"
$ws.Range("F19").Value = "'100"

# Grand total
$ws.Range("F29").Value = "'100"
